$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Multi-threading" topic questions appended below the existing
# Multi-threading rows (row 34 was the last used row, A1:C34 was the
# used range before this edit).
$questions = @(
    'How can we make sure main is the last thread to finish java program?',
    'How does thread communicate with each other?',
    'Why wait(), notify() and notifyAll() methods are in Object class?',
    'Why we call wait(), notify() and notifyAll() methods have to be called from synchronized method or block?',
    'Why sleep() and yeild() methods are static?',
    'Difference between interrupted() and isInterrupted() method.',
    'How can we achieve thread safty in java?',
    'Which is more preferred, synchronized method or block?',
    'What is ThreadLocal?'
)

$startRow = 35
for ($i = 0; $i -lt $questions.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = "Multi-threading"
    $ws.Cells.Item($row, 2).Value = $questions[$i]
}

# Match the author's final selection / scroll state from the diff
# (the sheet was scrolled so the frozen header row is followed by row 21,
# with the cursor left on B44, the first blank cell below the new rows).
try {
    $excel.ActiveWindow.ScrollRow = 21
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B44").Select()
